$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the one-hot "1" position within specific row pairs of columns,
# per the diff: (A3,B3), (E4,F4), (A6,B6), (E7,F7)

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 1

$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 0

$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0
